# "run all scenario sulawesi selatan"
# Updates the LDMProp_2006 scenario grid with the new run's percentages.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LDMProp_2006")

# --- Row 2 ---
$ws.Range("K2").Value = 0.05
$ws.Range("O2").Value = 0.999
$ws.Range("V2").Value = 0.03

# --- Rows 3-6 (same pattern) ---
$ws.Range("M3").Value = 0.2
$ws.Range("R3").Value = 0.05
$ws.Range("M4").Value = 0.2
$ws.Range("R4").Value = 0.05
$ws.Range("M5").Value = 0.2
$ws.Range("R5").Value = 0.05
$ws.Range("M6").Value = 0.2
$ws.Range("R6").Value = 0.05

# --- Rows 7-9 ---
$ws.Range("R7").Value = 0.05
$ws.Range("R8").Value = 0.05
$ws.Range("R9").Value = 0.05

# --- Row 10: N10 becomes a formula, R10 updated ---
$ws.Range("N10").Formula = "=1-0.701"
$ws.Range("R10").Value = 0.05

# --- Row 14 ---
$ws.Range("M14").Value = 0.2

# --- Row 15 ---
$ws.Range("J15").Value = 0.15
$ws.Range("L15").Value = 0.25

# --- Rows 16-17: Q formula 0.3/39 -> 0.2/39 (new shared formula group) ---
$ws.Range("Q16:Q17").Formula = "=0.2/39"

# --- Row 18 ---
$ws.Range("B18").Value = 0.2
$ws.Range("D18").Value = 0
$ws.Range("E18").Value = 0.3
$ws.Range("F18").Value = 0.2
$ws.Range("G18").Value = 0.2
$ws.Range("H18").Value = 0.9
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("V18").Value = 0

# --- Row 19 ---
$ws.Range("F19").Value = 0.2
$ws.Range("G19").Value = 0.2
$ws.Range("P19").Value = 0.55
$ws.Range("U19").Value = 0.3
$ws.Range("V19").Value = 0.01

# --- Rows 20-23 ---
$ws.Range("S20").Value = 0.05
$ws.Range("S21").Value = 0.05
$ws.Range("S22").Value = 0.05
$ws.Range("S23").Value = 0.05

# --- Row 24: standalone Q formula 0.3/39 -> 0.2/39 ---
$ws.Range("Q24").Formula = "=0.2/39"

# --- Rows 25-60: shared Q formula 0.3/39 -> 0.2/39 ---
$ws.Range("Q25:Q60").Formula = "=0.2/39"

# --- Row 48, 54, 58, 59 (X / C tweaks) ---
$ws.Range("X48").Value = 0
$ws.Range("X54").Value = 0.9
$ws.Range("X58").Value = 0.1
$ws.Range("C59").Value = 0.1

# --- Row 62: full scenario row rewrite ---
$ws.Range("B62").Value = 0.8
$ws.Range("C62").Value = 0.6
$ws.Range("D62").Value = 1
$ws.Range("E62").Value = 0.7
$ws.Range("F62").Value = 0.6
$ws.Range("G62").Value = 0.6
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0.85
$ws.Range("K62").Value = 0.95
$ws.Range("L62").Value = 0.75
$ws.Range("M62").Value = 0
$ws.Range("N62").Value = 0.001
$ws.Range("O62").Value = 0.001
$ws.Range("P62").Value = 0.45
$ws.Range("Q62").Value = 0.8
$ws.Range("R62").Value = 0.5
$ws.Range("S62").Value = 0.8
$ws.Range("T62").Value = 0
$ws.Range("U62").Value = 0.7
$ws.Range("V62").Value = 0.96
$ws.Range("W62").Value = 1
$ws.Range("X62").Value = 0

# --- Row 63: O63 becomes a SUM formula like the rest of the row ---
$ws.Range("O63").Formula = "=SUM(O2:O62)"

# --- Column widths (best achievable via the character-width COM property) ---
$ws.Columns.Item(1).ColumnWidth = 28.65
$ws.Columns.Item(7).ColumnWidth = 10.8
$ws.Columns.Item(11).ColumnWidth = 8.65
$ws.Columns.Item(13).ColumnWidth = 11.0

# --- Selection shown in the saved view ---
$ws.Range("B1:X62").Select()
